$wb = $excel.ActiveWorkbook

function Set-HitRow($ws, $row, $hits, $pct) {
    $ws.Range("B$row").Value = $hits
    $ws.Range("D$row").NumberFormat = "@"
    $ws.Range("D$row").Value = $pct
}

# --- Sheet "Total Hits" ---
$ws = $wb.Worksheets.Item("Total Hits")
Set-HitRow $ws 2 1100 "37.28%"
Set-HitRow $ws 3 2174 "36.83%"
Set-HitRow $ws 4 3266 "36.89%"
Set-HitRow $ws 5 4376 "37.07%"
Set-HitRow $ws 6 5466 "37.05%"

# --- Sheet "Hits_entity" ---
$ws = $wb.Worksheets.Item("Hits_entity")
Set-HitRow $ws 3 920 "31.72%"
Set-HitRow $ws 4 1369 "31.47%"
Set-HitRow $ws 5 1839 "31.71%"
Set-HitRow $ws 6 2302 "31.75%"

# --- Sheet "Hits_numerical" ---
$ws = $wb.Worksheets.Item("Hits_numerical")
Set-HitRow $ws 2 156 "23.82%"
Set-HitRow $ws 3 322 "24.58%"
Set-HitRow $ws 4 492 "25.04%"
Set-HitRow $ws 5 658 "25.11%"
Set-HitRow $ws 6 809 "24.70%"

# --- Sheet "Hits_date" ---
$ws = $wb.Worksheets.Item("Hits_date")
Set-HitRow $ws 3 220 "41.51%"
Set-HitRow $ws 4 328 "41.26%"
Set-HitRow $ws 5 439 "41.42%"
Set-HitRow $ws 6 562 "42.42%"
